$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Combo resettable fuse and RCD" row (row 15):
# type -> "KS 68 B16A LS-DI", alternative type -> "ABB DS202C", price per unit -> 136.8
$ws.Range("E15").Value = "ABB DS202C"
$ws.Range("D15").Value = "KS 68 B16A LS-DI"
$ws.Range("G15").Value = 136.8

# Move the active selection like in the final saved state
$ws.Range("D15").Select()
